# Update (Removed Auto Arima)
# Fill in the previously-blank "Amazon Mean Forecast" column (D) and
# refresh the P70/P80/P90 forecast columns (E/F/G) on the
# "Forecast Comparison" sheet now that the Auto-ARIMA model has been
# removed from the forecast ensemble.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Row -> D, E, F, G values
$data = @{
    2  = @(10, 11, 17, 27)
    3  = @(9,  9,  15, 25)
    4  = @(8,  8,  13, 23)
    5  = @(9,  9,  15, 25)
    6  = @(9,  9,  14, 25)
    7  = @(9,  9,  15, 25)
    8  = @(8,  9,  14, 25)
    9  = @(8,  7,  13, 25)
    10 = @(8,  8,  14, 24)
    11 = @(8,  8,  13, 24)
    12 = @(8,  7,  13, 24)
    13 = @(9,  8,  14, 26)
    14 = @(8,  7,  13, 24)
    15 = @(8,  6,  12, 24)
    16 = @(8,  7,  13, 24)
    17 = @(7,  6,  12, 23)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
}
